$d = $word.ActiveDocument

# --- Step 1: remove the stray "_GoBack" bookmark that currently sits inside
# the R9 paragraph (between "the " and "help of a map API"). Word will
# re-create it at the location of the most recent edit once we're done.
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # no-op if it doesn't exist for some reason
}

# --- Step 2: locate the "D2" domain-assumption paragraph that belongs to the
# G2 goal (the second "D2: ..." paragraph in the document, the one that
# follows "R12: ...") and remember its paragraph index.
$range = $d.Content
$range.Find.Execute("R12: If accident reports are provided by authorities the system should take that data into account when calculating the safety of a certain area", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$range.Collapse(0)
$range.Find.Execute("D2: User location included in report assumed to be the true unmodified location", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$range.Collapse(0)

$d2Paragraph = $range.Paragraphs(1)
$d2Index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $d2Paragraph.Range.Start) {
        $d2Index = $i
        break
    }
}

# --- Step 3: insert a brand-new paragraph right after that D2 paragraph and
# fill it in with the new domain assumption D9.
$range.InsertParagraphAfter()

$d9Index = $d2Index + 1
$d9Paragraph = $d.Paragraphs($d9Index)
$d9Paragraph.Range.Text = "D9: The communication of accident reports by the municipality is assumed to be proactive"

# --- Step 4: re-create the "_GoBack" bookmark at the end of the freshly
# typed D9 text, mirroring where Word leaves it after the last edit.
$d9End = $d.Paragraphs($d9Index).Range
$d9End.Collapse(0)
$d9End.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $d9End)
